$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in B4 (ratio_threshold_range / Min) from 1 to 0.9
$ws.Range("B4").Value = 0.9

# Move the active selection from B3 to B4, matching the saved selection state
[void]$ws.Range("B4").Select()
